$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 32 (ALC)
$ws.Range("H32").Value = 1078.4166
$ws.Range("I32").Value = 700
$ws.Range("J32").Value = 1154.1
$ws.Range("K32").Value = 700
$ws.Range("L32").Value = 1154.1
$ws.Range("M32").Value = -374
$ws.Range("N32").Value = -1806.1

# Row 51 (ALC)
$ws.Range("H51").Value = 2500
$ws.Range("I51").Value = 1933.3334
$ws.Range("K51").Value = 1933.3334
$ws.Range("M51").Value = -1449.3334

# Row 112 (ALC)
$ws.Range("H112").Value = 11905851
$ws.Range("I112").Value = 638.5714
$ws.Range("J112").Value = 14286893
$ws.Range("K112").Value = 1915.7142
$ws.Range("L112").Value = 42860679
$ws.Range("M112").Value = -807.7142000000001
$ws.Range("N112").Value = -42862895

# Row 137 (ALC)
$ws.Range("H137").Value = 1408.0392
$ws.Range("I137").Value = 829.5789
$ws.Range("J137").Value = 1751.5
$ws.Range("K137").Value = 2488.7367
$ws.Range("L137").Value = 5254.5
$ws.Range("M137").Value = 61.26330000000007
$ws.Range("N137").Value = -10354.5

# Row 139 (ALC)
$ws.Range("H139").Value = 0
$ws.Range("J139").Value = 0
$ws.Range("L139").Value = 0
$ws.Range("N139").ClearContents()

# Row 141 (ALC)
$ws.Range("H141").Value = 1943.7142
$ws.Range("I141").Value = 1943.7142
$ws.Range("K141").Value = 5831.142599999999
$ws.Range("M141").Value = -651.1425999999992

$ws = $wb.Worksheets.Item("ARM")
# Row 32 (ARM)
$ws.Range("H32").Value = 7692.94
$ws.Range("I32").Value = 8186.4116
$ws.Range("J32").Value = 6644.3125
$ws.Range("K32").Value = 8186.4116
$ws.Range("L32").Value = 6644.3125
$ws.Range("M32").Value = -7899.4116
$ws.Range("N32").Value = -7218.3125

# Row 61 (ARM)
$ws.Range("H61").Value = 3796.3022
$ws.Range("I61").Value = 4797.2256
$ws.Range("K61").Value = 4797.2256
$ws.Range("M61").Value = -4585.2256

# Row 74 (ARM)
$ws.Range("H74").Value = 5358.8276
$ws.Range("I74").Value = 7845.375
$ws.Range("J74").Value = 2298.4614
$ws.Range("K74").Value = 7845.375
$ws.Range("L74").Value = 2298.4614
$ws.Range("M74").Value = -6971.375
$ws.Range("N74").Value = -4046.4614

# Row 77 (ARM)
$ws.Range("H77").Value = 5358.8276
$ws.Range("I77").Value = 7845.375
$ws.Range("J77").Value = 2298.4614
$ws.Range("K77").Value = 39226.875
$ws.Range("L77").Value = 11492.307
$ws.Range("M77").Value = -34858.875
$ws.Range("N77").Value = -20228.307

# Row 136 (ARM)
$ws.Range("H136").Value = 3796.3022
$ws.Range("I136").Value = 4797.2256
$ws.Range("K136").Value = 14391.6768
$ws.Range("M136").Value = -11841.6768

$ws = $wb.Worksheets.Item("CRP")
# Row 31 (CRP)
$ws.Range("H31").Value = 2580.5059
$ws.Range("I31").Value = 1864.174
$ws.Range("J31").Value = 2846.242
$ws.Range("K31").Value = 1864.174
$ws.Range("L31").Value = 2846.242
$ws.Range("M31").Value = -1569.174
$ws.Range("N31").Value = -3436.242

# Row 34 (CRP)
$ws.Range("H34").Value = 2580.5059
$ws.Range("I34").Value = 1864.174
$ws.Range("J34").Value = 2846.242
$ws.Range("K34").Value = 1864.174
$ws.Range("L34").Value = 2846.242
$ws.Range("M34").Value = -1662.174
$ws.Range("N34").Value = -3250.242

# Row 38 (CRP)
$ws.Range("H38").Value = 6000
$ws.Range("I38").Value = 6000
$ws.Range("K38").Value = 6000
$ws.Range("M38").Value = -5623

# Row 39 (CRP)
$ws.Range("H39").Value = 6779.25
$ws.Range("I39").Value = 0
$ws.Range("J39").Value = 6779.25
$ws.Range("K39").Value = 0
$ws.Range("L39").Value = 6779.25
$ws.Range("M39").ClearContents()
$ws.Range("N39").Value = -7561.25

# Row 42 (CRP)
$ws.Range("H42").Value = 11285.333
$ws.Range("J42").Value = 13900
$ws.Range("L42").Value = 13900
$ws.Range("N42").Value = -15086

# Row 46 (CRP)
$ws.Range("H46").Value = 6000
$ws.Range("I46").Value = 6000
$ws.Range("K46").Value = 6000
$ws.Range("M46").Value = -5789

# Row 48 (CRP)
$ws.Range("H48").Value = 25000
$ws.Range("J48").Value = 25000
$ws.Range("L48").Value = 25000
$ws.Range("N48").Value = -25952

# Row 49 (CRP)
$ws.Range("H49").Value = 6779.25
$ws.Range("I49").Value = 0
$ws.Range("J49").Value = 6779.25
$ws.Range("K49").Value = 0
$ws.Range("L49").Value = 6779.25
$ws.Range("M49").ClearContents()
$ws.Range("N49").Value = -7143.25

# Row 50 (CRP)
$ws.Range("H50").Value = 19849.5
$ws.Range("J50").Value = 19849.5
$ws.Range("L50").Value = 19849.5
$ws.Range("N50").Value = -21099.5

# Row 54 (CRP)
$ws.Range("H54").Value = 0
$ws.Range("J54").Value = 0
$ws.Range("L54").Value = 0
$ws.Range("N54").ClearContents()

# Row 58 (CRP)
$ws.Range("H58").Value = 1954.625
$ws.Range("I58").Value = 1613.6666
$ws.Range("J58").Value = 2159.2
$ws.Range("K58").Value = 1613.6666
$ws.Range("L58").Value = 2159.2
$ws.Range("M58").Value = -1410.6666
$ws.Range("N58").Value = -2565.2

# Row 122 (CRP)
$ws.Range("H122").Value = 2980
$ws.Range("I122").Value = 2980
$ws.Range("K122").Value = 8940
$ws.Range("M122").Value = -6490

# Row 136 (CRP)
$ws.Range("H136").Value = 1954.625
$ws.Range("I136").Value = 1613.6666
$ws.Range("J136").Value = 2159.2
$ws.Range("K136").Value = 4840.9998
$ws.Range("L136").Value = 6477.599999999999
$ws.Range("M136").Value = -2290.9998
$ws.Range("N136").Value = -11577.6

$ws = $wb.Worksheets.Item("CUL")
# Row 68 (CUL)
$ws.Range("H68").Value = 2430.1099
$ws.Range("I68").Value = 2795.9805
$ws.Range("J68").Value = 1963.625
$ws.Range("K68").Value = 8387.941500000001
$ws.Range("L68").Value = 5890.875
$ws.Range("M68").Value = -7576.941500000001
$ws.Range("N68").Value = -7512.875

# Row 71 (CUL)
$ws.Range("H71").Value = 2430.1099
$ws.Range("I71").Value = 2795.9805
$ws.Range("J71").Value = 1963.625
$ws.Range("K71").Value = 25163.8245
$ws.Range("L71").Value = 17672.625
$ws.Range("M71").Value = -21107.8245
$ws.Range("N71").Value = -25784.625

# Row 107 (CUL)
$ws.Range("H107").Value = 1203.5209
$ws.Range("I107").Value = 303.66666
$ws.Range("J107").Value = 1411.1794
$ws.Range("K107").Value = 910.9999799999999
$ws.Range("L107").Value = 4233.5382
$ws.Range("M107").Value = 1009.00002
$ws.Range("N107").Value = -8073.5382

# Row 118 (CUL)
$ws.Range("H118").Value = 4940
$ws.Range("I118").Value = 0
$ws.Range("J118").Value = 4940
$ws.Range("K118").Value = 0
$ws.Range("L118").Value = 14820
$ws.Range("M118").ClearContents()
$ws.Range("N118").Value = -17306

# Row 131 (CUL)
$ws.Range("H131").Value = 11963213
$ws.Range("J131").Value = 12988419
$ws.Range("L131").Value = 38965257
$ws.Range("N131").Value = -38975337

$ws = $wb.Worksheets.Item("GSM")
# Row 132 (GSM)
$ws.Range("H132").Value = 3875.2432
$ws.Range("I132").Value = 4826.222
$ws.Range("K132").Value = 14478.666
$ws.Range("M132").Value = -11948.666

$ws = $wb.Worksheets.Item("WVR")
# Row 57 (WVR)
$ws.Range("H57").Value = 52400
$ws.Range("I57").Value = 52400
$ws.Range("K57").Value = 52400
$ws.Range("M57").Value = -51646

# Row 113 (WVR)
$ws.Range("H113").Value = 778.8570999999999
$ws.Range("I113").Value = 624.4138
$ws.Range("J113").Value = 1123.3846
$ws.Range("K113").Value = 1873.2414
$ws.Range("L113").Value = 3370.1538
$ws.Range("M113").Value = 296.7585999999999
$ws.Range("N113").Value = -7710.1538

# Row 136 (WVR)
$ws.Range("H136").Value = 1566.1852
$ws.Range("I136").Value = 879.875
$ws.Range("K136").Value = 2639.625
$ws.Range("M136").Value = -89.625
